$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "estado-de-la-informacion" column (J) moves from an iaest-dimension to an
# iaest-measure, so its type-row values change from dim/skos:Concept to
# medida/xsd:int.
$ws.Range("J2").Value = "iaest-measure:estado-de-la-informacion"
$ws.Range("J3").Value = "medida"
$ws.Range("J4").Value = "xsd:int"

# "municipio-nombre" column (M) moves from an iaest-measure to the
# sdmx-dimension:refArea dimension (like provincia-nombre/comarca-nombre),
# so its type-row values change from medida/xsd:int to dim/URI-Municipio.
$ws.Range("M2").Value = "sdmx-dimension:refArea"
$ws.Range("M3").Value = "dim"
$ws.Range("M4").Value = "URI-Municipio"

# The mapping file reference for estado-de-la-informacion is no longer
# needed now that it isn't a curated dimension with its own mapping sheet.
$ws.Range("J5").Clear()
